# "algoritmul lui lee pentru gasirea path-ului"
# Marks attendance in column L ("sapt. 11") for the students who presented
# Lee's algorithm (path finding) during week 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose value in column L ("sapt. 11") already existed and must be
# incremented from 1 to 2 (a second mark / extra point for that week).
$incrementRows = @(4, 5, 8, 11, 13, 19, 20, 25, 34, 37, 44, 46, 48, 49)
foreach ($r in $incrementRows) {
    $ws.Cells.Item($r, 12).Value = 2
}

# Rows that had no value in column L yet and now get a fresh mark of 1.
$newOneRows = @(7, 9, 23, 24, 26, 41, 45)
foreach ($r in $newOneRows) {
    $ws.Cells.Item($r, 12).Value = 1
}

# Row 18 had no value in column L yet and now gets a fresh mark of 2.
$ws.Cells.Item(18, 12).Value = 2

# Update the current selection to match the end state of the edit.
$ws.Range("M9").Select()
